$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New StatQuery text (replaces the old StatQuery column text for the Cases/Samples/Files rows)
$newStatQuery = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Giant Schnauzer']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@

# Update the StatQuery column (C) for the Cases (row 2), Samples (row 3) and Files (row 4) rows
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Column D ("dbExcel") no longer needs to be as wide - autofit it to the (unchanged) file-name content
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 44.3

# Move the active selection from B2 to B4
$ws.Range("B4").Select() | Out-Null

Write-Output "Applied ICDC breed testcase fix"
